$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'95.154.12"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.24%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'3.572.20"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.72%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  +0.10%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'235.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.93%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'652.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +2.29%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "'  -0.76%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.399"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.58%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("E9").Value = "'  +0.13%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("E10").Value = "'  -2.12%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'3.572.98"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.68%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("E12").Value = "'  +1.12%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'42.45"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.67%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("E14").Value = "'  +1.58%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'4.239.73"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.75%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'95.088.18"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -1.26%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'0.0000253"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.13%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'3.580.17"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.42%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("E19").Value = "'  -3.47%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'12.59"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -4.95%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'17.84"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -1.50%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'3.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.01%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'508.73"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.34%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'0.478"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -4.95%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'6.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +2.01%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'0.0000195"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -1.72%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'95.33"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.55%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("E28").Value = "'  +2.09%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'3.765.04"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.48%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("E30").Value = "'  -1.37%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("E31").Value = "'  -0.78%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("E32").Value = "'  -0.75%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("E33").Value = "'  +0.39%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("E34").Value = "'  -0.13%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D36").Value = "'31.70"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +4.24%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'1.68"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +12.93%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'8.60"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +9.42%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'0.558"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.88%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'583.57"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +1.22%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("E41").Value = "'  +0.03%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("E42").Value = "'  -0.90%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'0.905"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -2.26%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'1.81"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +3.40%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'2.30"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +5.18%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'5.74"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +1.41%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("B47").Value = "'EnergySwap"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'33.93"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +30.23%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("B48").Value = "'WhiteBITCoin"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'23.39"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -1.72%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("E49").Value = "'  -3.74%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'3.58"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.46%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'8.15"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.33%  "
$ws.Range("E51").Style = "Normal"
